# RAP4.xlsx fix: the RAP-generated html-templates ("PF_Interface" / navigation-menu
# editor rows) were removed from the "Identity Provider data" sheet together with
# their now-unused shared strings, the "Janitor"/"Advanced" role rows were swapped
# back into their original order, and the selection/scroll position was updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Identity Provider data")

# --- Swap the "Janitor" / "Advanced" rows back (34 <-> 35) -----------------
$ws.Range("A34").Value = "Advanced"
$ws.Range("B34").Value = "Advanced"
$ws.Range("A35").Value = "Janitor"
$ws.Range("B35").Value = "Janitor"

# --- Move the PF_NavMenuItem/Login/Logout/Anonymous block up into rows 37-41
#     (it used to live in rows 44-48) ---------------------------------------
$ws.Range("A37").Value = "[PF_NavMenuItem]"
$ws.Range("B37").Value = "isPartOf"
$ws.Range("C37").Value = "label"

$ws.Range("A38").Value = "PF_NavMenuItem"
$ws.Range("B38").Value = "PF_NavMenu"
$ws.Range("C38").Value = "PF_Label"

$ws.Range("A39").Value = "Login"
$ws.Range("B39").Value = "Login"
$ws.Range("C39").Value = "Login"

$ws.Range("A40").Value = "Logout"
$ws.Range("B40").Value = "Logout"
$ws.Range("C40").Value = "Logout"

$ws.Range("A41").Value = "Anonymous"
$ws.Range("B41").Value = "Anonymous"
$ws.Range("C41").Value = "Anonymous"

# Re-apply the correct cell formatting for the moved block:
#  - rows 37-38 use the "label" look (same format as A1, style index 1)
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A37:C38").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

#  - rows 39-41 use the plain look already used elsewhere in this table (style index 10)
$ws.Range("A33").Copy() | Out-Null
$ws.Range("A39:C41").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

# --- Open up a blank row (new row 42) before the existing blank row 43 -----
$ws.Rows("43").Insert(-4121) | Out-Null               # xlShiftDown

# Give the two blank rows (42 and 43) explicit formatting consistent with the
# rest of the blank rows in this block (style index 10, like row 36).
$ws.Range("A36:D36").Copy() | Out-Null
$ws.Range("A42:D43").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Remove the now-duplicated old rows (old 44-50, shifted to 44-51 by the
#     insert above) which used to hold the PF_Interface / navigation-menu
#     editor template rows that are no longer generated by RAP -------------
$ws.Range("A44:G51").EntireRow.Delete()

# --- Update selection / scroll position ------------------------------------
$ws.Range("B35").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1

Write-Output ("Final UsedRange: " + $ws.UsedRange.Address())
